$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Attacco di Hamas a Israele" / La Repubblica / Instagram -> num_commenti 93 -> 100
$ws.Range("E10").Value = 100

# Row 33: "La morte di Silvio Berlusconi" / Il Corriere Della Sera / Facebook -> num_commenti 99 -> 100
$ws.Range("E33").Value = 100

# Row 36 had a typo'd journal name and a bad comment count; fix it to match the
# "La Repubblica" pattern used by the surrounding rows.
$ws.Range("C36").Value = "La Repubblica"
$ws.Range("E36").Value = 100

# Row 37 social channel was Facebook (duplicate of row 36); change to Instagram.
$ws.Range("D37").Value = "Instagram"

# Row 38 social channel was Instagram (duplicate); change to YouTube.
$ws.Range("D38").Value = "YouTube"

# Former row 39 (La morte di Silvio Berlusconi / La Repubblica / YouTube / 100) is now
# redundant with the corrected row 38, so remove it entirely, shrinking the sheet.
$ws.Range("A39:E39").Delete()
